# Update the "LR" (loren) GitHub paths in column B: the GitHub organization
# folder was renamed from "MARIO Organization" to "SESAM".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "C:\Users\loren\Documents\GitHub\SESAM\GreenTechs\Database"
$ws.Range("B6").Value = "C:\Users\loren\Documents\GitHub\SESAM\GreenTechs\Add sectors"
$ws.Range("B7").Value = "C:\Users\loren\Documents\GitHub\SESAM\GreenTechs\Shocks"
$ws.Range("B8").Value = "C:\Users\loren\Documents\GitHub\SESAM\GreenTechs\Results"
$ws.Range("B9").Value = "C:\Users\loren\Documents\GitHub\SESAM\GreenTechs\Plots"
$ws.Range("B10").Value = "C:\Users\loren\Documents\GitHub\SESAM\GreenTechs\Shocks\ShockMaster.xlsx"

# Reflect the updated selection / active cell after editing column B.
$ws.Range("B11").Select()
